$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rebuild Tasks")

# 1. Update text of two existing task rows.
$ws.Range("A9").Value = " -- Create User Account / Login / Logout"
$ws.Range("A12").Value = " -- Home Page & View FAQs"

# 2. Add the new "Follow-up Schedule" block in columns G:J starting at row 4.
$ws.Range("G4").Value = "Rebuild Task Update Schedule:"
$ws.Range("G4").Style = "Normal"
$ws.Range("G4").Font.Bold = $true
$ws.Range("G4").Font.Underline = $true
$ws.Range("G4").Font.Size = 14
$ws.Range("G4:J4").Interior.Pattern = -4124
$ws.Range("G4:J4").Interior.PatternColorIndex = -4105
$ws.Range("G4:J4").Interior.Color = 65535

$dateRows = @(5, 6, 7, 9, 10, 11, 12, 13, 14)
$dates = @(43088, 43109, 43116, 43123, 43130, 43137, 43144, 43151, 43158)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $cell = $ws.Range("G" + $dateRows[$i])
    $cell.Value = [DateTime]::FromOADate($dates[$i])
    $cell.NumberFormat = "m/d/yyyy"
    $cell.Font.Bold = $true
    $cell.Font.Size = 12
}

# Row 8 is hidden; its G cell exists but stays blank (still gets the bold font).
$ws.Range("G8").Font.Bold = $true
$ws.Range("G8").Font.Size = 12

$ws.Columns.Item(6).ColumnWidth = 5.7109375
$ws.Columns.Item(7).ColumnWidth = 12.140625

$ws.Range("I13").Select()
